# Update cryptocurrency price/volume data per the latest scrape.
# Row 47/48 also swap Coin/Link (Aave <-> Frax) in addition to new D/E values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '25.896.14'
$ws.Range("E2").Value2 = '  +0.24%  '

$ws.Range("D3").Value2 = '1.641.02'
$ws.Range("E3").Value2 = '  +0.44%  '

$ws.Range("D4").Value2 = '''1.009'
$ws.Range("E4").Value2 = '  +0.53%  '

$ws.Range("D5").Value2 = '''215.08'
$ws.Range("E5").Value2 = '  +0.08%  '

$ws.Range("D6").Value2 = '''0.5059'
$ws.Range("E6").Value2 = '  +0.88%  '

$ws.Range("D7").Value2 = '''1.006'
$ws.Range("E7").Value2 = '  +0.23%  '

$ws.Range("D8").Value2 = '''0.2575'
$ws.Range("E8").Value2 = '  +0.35%  '

$ws.Range("D9").Value2 = '''0.06408'
$ws.Range("E9").Value2 = '  -0.01%  '

$ws.Range("D10").Value2 = '''19.68'
$ws.Range("E10").Value2 = '  +0.73%  '

$ws.Range("D11").Value2 = '''0.07794'
$ws.Range("E11").Value2 = '  +1.46%  '

$ws.Range("D12").Value2 = '''4.293'
$ws.Range("E12").Value2 = '  +1.62%  '

$ws.Range("D13").Value2 = '1.630.05'
$ws.Range("E13").Value2 = '  -0.27%  '

$ws.Range("D14").Value2 = '''0.5450'
$ws.Range("E14").Value2 = '  +0.02%  '

$ws.Range("D15").Value2 = '0.0₅7884'
$ws.Range("E15").Value2 = '  -0.37%  '

$ws.Range("D16").Value2 = '''65.12'
$ws.Range("E16").Value2 = '  +2.69%  '

$ws.Range("D17").Value2 = '25.969.54'
$ws.Range("E17").Value2 = '  +0.45%  '

$ws.Range("D18").Value2 = '''1.007'
$ws.Range("E18").Value2 = '  +0.38%  '

$ws.Range("D19").Value2 = '''198.15'
$ws.Range("E19").Value2 = '  -2.31%  '

$ws.Range("D20").Value2 = '''4.404'
$ws.Range("E20").Value2 = '  +2.47%  '

$ws.Range("D21").Value2 = '''9.975'
$ws.Range("E21").Value2 = '  +0.48%  '

$ws.Range("D22").Value2 = '''6.049'
$ws.Range("E22").Value2 = '  +1.32%  '

$ws.Range("D23").Value2 = '''1.007'
$ws.Range("E23").Value2 = '  +0.26%  '

$ws.Range("D24").Value2 = '''1.868'
$ws.Range("E24").Value2 = '  -3.56%  '

$ws.Range("D25").Value2 = '''140.45'
$ws.Range("E25").Value2 = '  -0.34%  '

$ws.Range("D26").Value2 = '''0.1145'
$ws.Range("E26").Value2 = '  +0.21%  '

$ws.Range("D27").Value2 = '''6.876'
$ws.Range("E27").Value2 = '  +2.82%  '

$ws.Range("E28").Value2 = '  +0.29%  '

$ws.Range("D29").Value2 = '''1.240'
$ws.Range("E29").Value2 = '  +0.27%  '

$ws.Range("D30").Value2 = '''0.05018'
$ws.Range("E30").Value2 = '  +0.96%  '

$ws.Range("D31").Value2 = '''3.268'
$ws.Range("E31").Value2 = '  -0.11%  '

$ws.Range("D32").Value2 = '''3.199'
$ws.Range("E32").Value2 = '  +0.67%  '

$ws.Range("D33").Value2 = '''1.539'
$ws.Range("E33").Value2 = '  +0.76%  '

$ws.Range("D34").Value2 = '''2.367'
$ws.Range("E34").Value2 = '  +0.54%  '

$ws.Range("D35").Value2 = '''0.8947'
$ws.Range("E35").Value2 = '  +0.56%  '

$ws.Range("D36").Value2 = '''2.599'
$ws.Range("E36").Value2 = '  -0.96%  '

$ws.Range("D37").Value2 = '1.134.33'
$ws.Range("E37").Value2 = '  -3.15%  '

$ws.Range("D38").Value2 = '''0.5537'
$ws.Range("E38").Value2 = '  -0.38%  '

$ws.Range("D39").Value2 = '''0.01559'
$ws.Range("E39").Value2 = '  +0.41%  '

$ws.Range("D40").Value2 = '''1.006'
$ws.Range("E40").Value2 = '  +0.34%  '

$ws.Range("D41").Value2 = '''5.686'
$ws.Range("E41").Value2 = '  +1.06%  '

$ws.Range("D42").Value2 = '''0.8153'
$ws.Range("E42").Value2 = '  +1.73%  '

$ws.Range("D43").Value2 = '''99.58'
$ws.Range("E43").Value2 = '  +0.51%  '

$ws.Range("E44").Value2 = '  +6.75%  '

$ws.Range("D45").Value2 = '1.777.30'
$ws.Range("E45").Value2 = '  +0.35%  '

$ws.Range("D46").Value2 = '''0.4545'
$ws.Range("E46").Value2 = '  +0.73%  '

$ws.Range("B47").Value2 = 'Aave'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value2 = '''55.28'
$ws.Range("E47").Value2 = '  +1.05%  '

$ws.Range("B48").Value2 = 'Frax'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value2 = '''1.005'
$ws.Range("E48").Value2 = '  +0.11%  '

$ws.Range("D49").Value2 = '''0.05091'
$ws.Range("E49").Value2 = '  +1.12%  '

$ws.Range("E50").Value2 = '  +0.33%  '

$ws.Range("D51").Value2 = '''0.09564'
$ws.Range("E51").Value2 = '  +3.32%  '
